$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.615385890007019
$ws.Range("B1").Value = 2.831643342971802
$ws.Range("C1").Value = 3.360304832458496
$ws.Range("D1").Value = 3.756871223449707
$ws.Range("E1").Value = 1.042465329170227
